$d = $word.ActiveDocument

function Escape-Xml($t) {
    return $t.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
}

function Build-Run($text, $preserve) {
    $t = Escape-Xml($text)
    if ($preserve) {
        return '<w:r><w:t xml:space="preserve">' + $t + '</w:t></w:r>'
    } else {
        return '<w:r><w:t>' + $t + '</w:t></w:r>'
    }
}

function Replace-WithRuns($OldText, $NewRunTexts, $PreserveSpace, $TrailingXml) {
    $full = $d.Content
    $full.Find.Execute($OldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

    $s = $full.Start
    $e = $full.End
    $rng = $d.Range($s, $e)

    $runsXml = ""
    for ($i = 0; $i -lt $NewRunTexts.Length; $i++) {
        $runsXml += Build-Run $NewRunTexts[$i] $PreserveSpace[$i]
    }
    $runsXml += $TrailingXml

    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $runsXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $rng.InsertXML($xml)
}

# 1) "Note: Sampler feedback ... a Scarlett only sample." (whole, single-run paragraph)
Replace-WithRuns `
    "Note: Sampler feedback is not supported on Xbox One, therefore this is a Scarlett only sample." `
    @(
        "Note: Sampler feedback is not supported on Xbox One, therefore this is ",
        "an Xbox Series X|S ",
        "only sample."
    ) `
    @($true, $true, $false) `
    ""

# 2) "This sample only supports Scarlett, so the active solution platform will be " is the
#    paragraph's FIRST run; it is immediately followed by a <w:proofErr/> marker, and this
#    runtime's InsertXML mis-orders any non-run sibling that trails a replacement touching the
#    paragraph's first child. Work around it by matching (and re-emitting) the whole paragraph,
#    so nothing is left trailing the inserted runs.
$run_gaming = Build-Run "Gaming.Xbox.Scarlett.x" $false
$run_64 = Build-Run "64" $false
$trail2 = '<w:proofErr w:type="gramStart"/>' + $run_gaming + '<w:proofErr w:type="gramEnd"/>' + $run_64

Replace-WithRuns `
    "This sample only supports Scarlett, so the active solution platform will be Gaming.Xbox.Scarlett.x64" `
    @(
        "This sample only supports ",
        "Xbox Series X|S",
        ", so the active solution platform will be "
    ) `
    @($true, $false, $true) `
    $trail2

# 3) ". This sample creates a readback texture ... On Scarlett, the values ... 5.3 fixed point."
#    (this run is the paragraph's LAST child, so no trailing-sibling reordering risk)
Replace-WithRuns `
    ". This sample creates a readback texture which is used for the readback. On Scarlett, the values in the feedback map are 5.3 fixed point." `
    @(
        ". This sample creates a readback texture which is used for the readback. On ",
        "Xbox Series X|S",
        ", the values in the feedback map are 5.3 fixed point."
    ) `
    @($true, $false, $false) `
    ""
